$wb = $excel.ActiveWorkbook

# ALC row 9: Distill, My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 236.66667
$ws.Range("I9").Value = 220.90909
$ws.Range("J9").Value = 280
$ws.Range("K9").Value = 220.90909
$ws.Range("L9").Value = 280
$ws.Range("M9").Value = -51.90908999999999
$ws.Range("N9").Value = -618

# ALC row 19: Unbreak My Heart
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1271.7142
$ws.Range("I19").Value = 1131.8
$ws.Range("J19").Value = 1349.4445
$ws.Range("K19").Value = 1131.8
$ws.Range("L19").Value = 1349.4445
$ws.Range("M19").Value = -956.8
$ws.Range("N19").Value = -1699.4445

# ALC row 129: Practical Command
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 916.5700000000001
$ws.Range("I129").Value = 554.2143
$ws.Range("J129").Value = 975.55817
$ws.Range("K129").Value = 1662.6429
$ws.Range("L129").Value = 2926.67451
$ws.Range("M129").Value = 3337.3571
$ws.Range("N129").Value = -12926.67451

# ALC row 135: For Tired Minds
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1619.2632
$ws.Range("I135").Value = 582
$ws.Range("J135").Value = 5509
$ws.Range("K135").Value = 5238
$ws.Range("L135").Value = 49581
$ws.Range("M135").Value = -2703

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1553.9048
$ws.Range("I137").Value = 1581.6
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 4744.799999999999
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -2194.799999999999

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26646.521
$ws.Range("I32").Value = 5535.6567
$ws.Range("J32").Value = 380253.5
$ws.Range("K32").Value = 5535.6567
$ws.Range("L32").Value = 380253.5
$ws.Range("M32").Value = -5248.6567

# ARM row 45: Hollow Hallmarks
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 101479.5
$ws.Range("I45").Value = 126349.75
$ws.Range("J45").Value = 1998.5
$ws.Range("K45").Value = 126349.75
$ws.Range("L45").Value = 1998.5
$ws.Range("M45").Value = -125972.75

# ARM row 80: A Squire to Inspire
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 27453.2
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 27453.2
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 27453.2
$ws.Range("N80").Value = -29449.2

# ARM row 83: All's Fair in Highborn Assassination (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 27453.2
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 27453.2
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 82359.60000000001
$ws.Range("N83").Value = -92343.60000000001

# ARM row 122: Haste for High Durium
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1940.2778
$ws.Range("I122").Value = 1763
$ws.Range("J122").Value = 2826.6667
$ws.Range("K122").Value = 5289
$ws.Range("L122").Value = 8480.000100000001
$ws.Range("M122").Value = -2839
$ws.Range("N122").Value = -13380.0001

# BSM row 22: Riveting Run
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 426.91666
$ws.Range("I22").Value = 447.72726
$ws.Range("J22").Value = 198
$ws.Range("K22").Value = 447.72726
$ws.Range("L22").Value = 198
$ws.Range("M22").Value = -274.72726
$ws.Range("N22").Value = -544

# BSM row 29: Powderpost Derby
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 2258
$ws.Range("I29").Value = 2258
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2258
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1969
$ws.Range("N29").ClearContents()

# BSM row 99: Meddle in Metal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1555.6666
$ws.Range("I99").Value = 1187.1428
$ws.Range("J99").Value = 1667.826
$ws.Range("K99").Value = 1187.1428
$ws.Range("L99").Value = 1667.826
$ws.Range("M99").Value = 310.8571999999999

# CRP row 22: Driving Up the Wall
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 317.375
$ws.Range("I22").Value = 219.85715
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 219.85715
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 130.14285

# CRP row 69: Landing the Big One
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 21045.4
$ws.Range("I69").Value = 18327.334
$ws.Range("J69").Value = 25122.5
$ws.Range("K69").Value = 18327.334
$ws.Range("L69").Value = 25122.5
$ws.Range("M69").Value = -17578.334
$ws.Range("N69").Value = -26620.5

# CRP row 72: Fishing for Profits (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 21045.4
$ws.Range("I72").Value = 18327.334
$ws.Range("J72").Value = 25122.5
$ws.Range("K72").Value = 54982.00199999999
$ws.Range("L72").Value = 75367.5
$ws.Range("M72").Value = -51238.00199999999
$ws.Range("N72").Value = -82855.5

# CRP row 94: Beech, Please
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1256.1765
$ws.Range("I94").Value = 1328
$ws.Range("J94").Value = 1234.0769
$ws.Range("K94").Value = 1328
$ws.Range("L94").Value = 1234.0769
$ws.Range("M94").Value = -877
$ws.Range("N94").Value = -2136.0769

# CRP row 99: O Pine
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9979
$ws.Range("I99").Value = 3369.2222
$ws.Range("J99").Value = 21876.6
$ws.Range("K99").Value = 3369.2222
$ws.Range("L99").Value = 21876.6
$ws.Range("M99").Value = -1871.2222

# CRP row 126: A Better Conductor
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9979
$ws.Range("I126").Value = 3369.2222
$ws.Range("J126").Value = 21876.6
$ws.Range("K126").Value = 10107.6666
$ws.Range("L126").Value = 65629.79999999999
$ws.Range("M126").Value = -7637.6666

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5135.9546
$ws.Range("I132").Value = 4843.375
$ws.Range("J132").Value = 5916.1665
$ws.Range("K132").Value = 14530.125
$ws.Range("L132").Value = 17748.4995
$ws.Range("M132").Value = -12000.125
$ws.Range("N132").Value = -22808.4995

# CUL row 8: Whip It
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 180.66667
$ws.Range("I8").Value = 180.66667
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 542.00001
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -403.00001

# CUL row 113: Can't Eat Just One
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 996.0952
$ws.Range("I113").Value = 1892.4286
$ws.Range("J113").Value = 547.9286
$ws.Range("K113").Value = 5677.2858
$ws.Range("L113").Value = 1643.7858
$ws.Range("M113").Value = -3507.2858
$ws.Range("N113").Value = -5983.7858

# GSM row 102: Put the Metal to the Peddle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 242356.8
$ws.Range("I102").Value = 1659.3846
$ws.Range("J102").Value = 503112.34
$ws.Range("K102").Value = 1659.3846
$ws.Range("L102").Value = 503112.34
$ws.Range("M102").Value = -37.38460000000009
$ws.Range("N102").Value = -506356.34

# LTW row 22: Skin off Their Backs
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1784.8572
$ws.Range("I22").Value = 1823.909
$ws.Range("J22").Value = 1741.9
$ws.Range("K22").Value = 1823.909
$ws.Range("L22").Value = 1741.9
$ws.Range("M22").Value = -1528.909

# LTW row 27: Fire and Hide
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1784.8572
$ws.Range("I27").Value = 1823.909
$ws.Range("J27").Value = 1741.9
$ws.Range("K27").Value = 1823.909
$ws.Range("L27").Value = 1741.9
$ws.Range("M27").Value = -1716.909

# LTW row 46: Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 920564.0600000001
$ws.Range("I46").Value = 397.8
$ws.Range("J46").Value = 1687369.4
$ws.Range("K46").Value = 397.8
$ws.Range("L46").Value = 1687369.4
$ws.Range("M46").Value = -209.8
$ws.Range("N46").Value = -1687745.4

# LTW row 93: Hide to Go Seek
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2325.647
$ws.Range("I93").Value = 2324
$ws.Range("J93").Value = 2331
$ws.Range("K93").Value = 2324
$ws.Range("L93").Value = 2331
$ws.Range("M93").Value = -1076

# LTW row 130: Generous Soles
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 39434.75
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 39434.75
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 39434.75
$ws.Range("N130").Value = -49474.75

# WVR row 43: Walk Softly and Carry a Big Halberd
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 100000
$ws.Range("I43").Value = 100000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 100000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -99851
$ws.Range("N43").ClearContents()

# WVR row 81: Where the Dragonflies, the Net Catches
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 500349.75
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 333799.66
$ws.Range("K81").Value = 2000000
$ws.Range("L81").Value = 667599.3199999999
$ws.Range("M81").Value = -1998939
$ws.Range("N81").Value = -669721.3199999999

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 500349.75
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 333799.66
$ws.Range("K84").Value = 10000000
$ws.Range("L84").Value = 3337996.6
$ws.Range("M84").Value = -9994696
$ws.Range("N84").Value = -3348604.6

# WVR row 116: All-purpose Overgarments
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 39660
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 39660
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 39660
$ws.Range("N116").Value = -48838
